# Update the "想去人数" (F column) figures for a handful of events on both
# the "展览" and "全部类型" worksheets, reflecting a refreshed data pull.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row number -> new F-column value
$updates = @{
    3  = 13797
    7  = 270
    10 = 85
    14 = 447
    15 = 5755
    19 = 84
    22 = 227
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
